$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily FX rows to append after the existing last row (906).
# Columns: datetime (Excel serial date), open, high, low, close, volume
$newRows = @(
    @(45461, 0.13778, 0.13782, 0.13778, 0.13782, 0),
    @(45462, 0.13782, 0.13783, 0.13775, 0.13775, 0),
    @(45463, 0.13775, 0.13775, 0.13769, 0.13769, 0),
    @(45464, 0.1377, 0.13771, 0.13767, 0.13767, 0),
    @(45467, 0.13768, 0.13776, 0.13768, 0.13773, 0),
    @(45468, 0.13772, 0.13775, 0.13765, 0.13765, 0),
    @(45469, 0.13765, 0.13766, 0.13761, 0.13761, 0),
    @(45470, 0.13761, 0.13762, 0.13753, 0.13753, 0),
    @(45471, 0.13755, 0.13766, 0.13755, 0.13756, 0),
    @(45474, 0.13756, 0.13757, 0.13754, 0.13754, 0),
    @(45475, 0.13754, 0.13755, 0.13749, 0.13749, 0),
    @(45476, 0.13749, 0.13755, 0.13749, 0.13751, 0),
    @(45477, 0.13751, 0.13757, 0.13751, 0.13756, 0),
    @(45478, 0.13757, 0.13761, 0.13753, 0.13753, 0),
    @(45481, 0.13755, 0.13756, 0.13754, 0.13754, 0),
    @(45482, 0.13754, 0.13755, 0.13746, 0.13746, 0),
    @(45483, 0.13746, 0.13748, 0.1374, 0.1374, 0),
    @(45484, 0.13741, 0.13788, 0.13741, 0.13774, 0),
    @(45485, 0.13774, 0.13791, 0.13769, 0.13787, 0),
    @(45488, 0.13789, 0.13789, 0.13771, 0.13772, 0),
    @(45489, 0.13772, 0.13837, 0.13754, 0.13754, 0),
    @(45490, 0.13754, 0.13885, 0.13754, 0.13769, 0),
    @(45491, 0.13769, 0.13781, 0.13769, 0.13769, 0),
    @(45492, 0.1377, 0.13772, 0.13751, 0.13751, 0),
    @(45495, 0.13752, 0.13755, 0.13744, 0.13744, 0),
    @(45496, 0.13744, 0.13746, 0.13742, 0.13742, 0),
    @(45497, 0.13742, 0.13766, 0.13742, 0.13765, 0),
    @(45498, 0.13765, 0.13865, 0.13765, 0.13825, 0),
    @(45499, 0.13825, 0.13825, 0.13789, 0.13789, 0),
    @(45502, 0.13789, 0.1379, 0.13771, 0.13771, 0),
    @(45503, 0.13771, 0.13792, 0.13769, 0.13787, 0),
    @(45504, 0.13787, 0.13855, 0.13787, 0.13848, 0),
    @(45505, 0.13846, 0.13864, 0.138, 0.138, 0),
    @(45506, 0.13802, 0.1397, 0.13802, 0.13962, 0),
    @(45509, 0.13962, 0.14135, 0.13962, 0.14021, 0),
    @(45510, 0.14021, 0.1406, 0.1397, 0.1397, 0),
    @(45511, 0.1397, 0.14024, 0.13919, 0.13931, 0),
    @(45512, 0.13931, 0.13969, 0.13931, 0.13934, 0),
    @(45513, 0.13934, 0.13963, 0.13934, 0.13951, 0),
    @(45516, 0.1395, 0.13994, 0.13927, 0.13932, 0),
    @(45517, 0.13934, 0.14, 0.13927, 0.13972, 0),
    @(45518, 0.14048, 0.14051, 0.13983, 0.14008, 0),
    @(45519, 0.14007, 0.14035, 0.13964, 0.13978, 0)
)

$lastExistingRow = 906
$startRow = $lastExistingRow + 1
$endRow = $startRow + $newRows.Count - 1

# Column A carries a date-formatted style (same as the rest of the date
# column); copy that formatting down onto the new rows before writing values.
$ws.Range("A$lastExistingRow").Copy($ws.Range("A${startRow}:A${endRow}"))

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
